# 189. Rotate Array to right by K
# Add a new row (17) to the Leetcode tracker sheet for this question.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Problem number
$ws.Cells.Item(17, 1).Value = 189

# Title (new shared string "Rotate Array to right by K")
$ws.Cells.Item(17, 2).Value = "Rotate Array to right by K"
$ws.Cells.Item(17, 2).Style = "Normal"

# Language used
$ws.Cells.Item(17, 3).Value = "Java"

# Date solved (stored as a date-formatted serial number, matching the
# style used by the rows above it)
$ws.Cells.Item(17, 4).Value = 44992
$ws.Cells.Item(17, 4).NumberFormat = "d-mmm-yy"

# Move the active selection to the newly added cell, as happens after
# typing the last entry in the row.
$ws.Range("D17").Select()
